$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.137.87'
$ws.Range("E2").Value = '  +0.08%  '

$ws.Range("D3").Value = '1.798.66'
$ws.Range("E3").Value = '  +2.40%  '

$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").Value = '''336.02'
$ws.Range("E5").Value = '  -0.64%  '

$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("D7").Value = '''0.4616'
$ws.Range("E7").Value = '  +22.19%  '

$ws.Range("D8").Value = '''0.3702'
$ws.Range("E8").Value = '  +10.20%  '

$ws.Range("D9").Value = '''45.16'
$ws.Range("E9").Value = '  -0.37%  '

$ws.Range("D10").Value = '''0.07658'
$ws.Range("E10").Value = '  +6.27%  '

$ws.Range("D11").Value = '''1.148'
$ws.Range("E11").Value = '  +2.62%  '

$ws.Range("E12").Value = '  -0.01%  '

$ws.Range("D13").Value = '''1.002'
$ws.Range("E13").Value = '  +0.20%  '

$ws.Range("D14").Value = '''6.361'
$ws.Range("E14").Value = '  +3.33%  '

$ws.Range("E15").Value = '  +3.51%  '

$ws.Range("D16").Value = '1.796.36'
$ws.Range("E16").Value = '  +2.35%  '

$ws.Range("E17").Value = '  +3.46%  '

$ws.Range("E18").Value = '  +1.83%  '

$ws.Range("D19").Value = '''82.74'
$ws.Range("E19").Value = '  +2.80%  '

$ws.Range("E20").Value = '  +0.22%  '

$ws.Range("D21").Value = '''17.42'
$ws.Range("E21").Value = '  +3.18%  '

$ws.Range("D22").Value = '''6.412'
$ws.Range("E22").Value = '  +2.63%  '

$ws.Range("D23").Value = '28.141.22'
$ws.Range("E23").Value = '  +0.15%  '

$ws.Range("D24").Value = '''11.91'
$ws.Range("E24").Value = '  +2.08%  '

$ws.Range("D25").Value = '''2.413'
$ws.Range("E25").Value = '  +0.72%  '

$ws.Range("D26").Value = '''20.74'
$ws.Range("E26").Value = '  +4.59%  '

$ws.Range("D27").Value = '''2.381'
$ws.Range("E27").Value = '  +2.61%  '

$ws.Range("D28").Value = '''152.02'
$ws.Range("E28").Value = '  -0.89%  '

$ws.Range("D29").Value = '2.002.54'
$ws.Range("E29").Value = '  +2.40%  '

$ws.Range("D30").Value = '''134.12'
$ws.Range("E30").Value = '  +1.83%  '

$ws.Range("E31").Value = '  +1.43%  '

$ws.Range("D32").Value = '''4.050'
$ws.Range("E32").Value = '  +0.65%  '

$ws.Range("D33").Value = '''0.09610'
$ws.Range("E33").Value = '  +11.25%  '

$ws.Range("D34").Value = '''5.886'
$ws.Range("E34").Value = '  +1.67%  '

$ws.Range("D35").Value = '''0.02380'
$ws.Range("E35").Value = '  +2.44%  '

$ws.Range("D36").Value = '''0.2228'
$ws.Range("E36").Value = '  +5.63%  '

$ws.Range("E37").Value = '  -0.23%  '

$ws.Range("D38").Value = '''5.289'
$ws.Range("E38").Value = '  +2.60%  '

$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = '''0.6727'
$ws.Range("E39").Value = '  +0.80%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '''0.06357'
$ws.Range("E40").Value = '  +2.47%  '

$ws.Range("D41").Value = '''1.518'
$ws.Range("E41").Value = '  +4.82%  '

$ws.Range("D42").Value = '''1.237'
$ws.Range("E42").Value = '  +1.58%  '

$ws.Range("D43").Value = '''8.094'
$ws.Range("E43").Value = '  +0.73%  '

$ws.Range("D44").Value = '''14.17'
$ws.Range("E44").Value = '  +4.12%  '

$ws.Range("E45").Value = '  +0.10%  '

$ws.Range("D46").Value = '''0.6166'
$ws.Range("E46").Value = '  +1.92%  '

$ws.Range("D47").Value = '''3.843'
$ws.Range("E47").Value = '  +0.15%  '

$ws.Range("D48").Value = '''130.45'
$ws.Range("E48").Value = '  +1.35%  '

$ws.Range("D49").Value = '''2.060'
$ws.Range("E49").Value = '  +2.04%  '

$ws.Range("D50").Value = '''1.189'
$ws.Range("E50").Value = '  +1.48%  '

$ws.Range("D51").Value = '''0.07133'
$ws.Range("E51").Value = '  -0.20%  '
